# Applies the "first draft of comments" copy-edit pass to the Background /
# Issue / Methodology sections of the polyploidy proposal.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $r = $d.Content
    $ok = $r.Find.Execute($find, $true, $false, $false, $false, $false, `
                           $true, 1, $false, $replace, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $find"
    }
}

# --- Background, paragraph 1 (word choice) ---------------------------------
Replace-Text "studies on the extinction rates of polyploids propose that" `
             "studies on the extinction rates of polyploids suggest that"

# --- Background, paragraph 1 (typo + spacing + punctuation) ----------------
Replace-Text "critical to prediciting polyploidisation" `
             "critical to predicting polyploidisation"

Replace-Text "critical to predicting polyploidisation of plant species. Reproductive disadvantages" `
             "critical to predicting polyploidisation of plant species.  Reproductive disadvantages"

Replace-Text "Polyploids are frequently linked with distinct traits such as;" `
             "Polyploids are frequently linked with distinct traits such as"

Replace-Text "”$([char]32)which increase plant organs; reversal of selfing inhibition; enhanced capabilities for buffering of deleterious mutation (due to increased heterozygosity); and hybrid vigour (heterosis)" `
             "”, which include increases in plant organs, reversal of selfing inhibition, enhanced capabilities for buffering of deleterious mutation (due to increased heterozygosity), and hybrid vigour (heterosis)"

Replace-Text ". All of which are thought to counteract the reproductive disadvantages of polyploidy and instead make this mutation key to the invasive and adaptive potential of plants." `
             ". These traits are thought to overcome the reproductive disadvantages of polyploidy and instead make this mutation key to the invasive and adaptive potential of plants."

# --- Background, paragraph 2 (slash -> "and") -------------------------------
Replace-Text "autopolyploid/allopolyploid" "autopolyploid and allopolyploid"

# trailing space added at the end of the autopolyploid/allopolyploid paragraph
$p = $d.Paragraphs.Item(9)
$end = $d.Range($p.Range.Start, $p.Range.End - 1)
$end.Collapse(0)
$end.InsertAfter(" ")

# --- Issue paragraph (en-dash spacing) --------------------------------------
Replace-Text "Due–in part–to lack of predictive power" `
             "Due – in part – to lack of predictive power"

# trailing space added at the end of the Issue paragraph
$p = $d.Paragraphs.Item(11)
$end = $d.Range($p.Range.Start, $p.Range.End - 1)
$end.Collapse(0)
$end.InsertAfter(" ")

# --- Methodology paragraph (typo) -------------------------------------------
Replace-Text "how environmental disturbance effects polyploid distribution" `
             "how environmental disturbance affects polyploid distribution"

# trailing space added at the end of the Methodology paragraph
$p = $d.Paragraphs.Item(15)
$end = $d.Range($p.Range.Start, $p.Range.End - 1)
$end.Collapse(0)
$end.InsertAfter(" ")

# --- trailing space after "Is it suitable for CRAN?" ------------------------
$p = $d.Paragraphs.Item(18)
$end = $d.Range($p.Range.Start, $p.Range.End - 1)
$end.Collapse(0)
$end.InsertAfter(" ")

# --- trailing space after the Carnegie/12-weeks remark ----------------------
$p = $d.Paragraphs.Item(20)
$end = $d.Range($p.Range.Start, $p.Range.End - 1)
$end.Collapse(0)
$end.InsertAfter(" ")
